$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C header (13-01-2023), copy style/format from B1 (bold, centered, bordered)
$ws.Range("C1").Value = "13-01-2023"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Reorder rows: Fima Acciones, Fima PB Acciones, IEB Value now come first (rows 2-4),
# then avg (row 5), then total (row 6). Add column C values for each row.

# Row 2: Fima Acciones
$ws.Range("A2").Value = "Fima Acciones"
$ws.Range("B2").Value = 19546.74
$ws.Range("C2").Value = 19424.28

# Row 3: Fima PB Acciones
$ws.Range("A3").Value = "Fima PB Acciones"
$ws.Range("B3").Value = 10325.78
$ws.Range("C3").Value = 10368.85

# Row 4: IEB Value
$ws.Range("A4").Value = "IEB Value"
$ws.Range("B4").Value = 3846.39
$ws.Range("C4").Value = 3850.1

# Row 5: avg
$ws.Range("A5").Value = "avg"
$ws.Range("B5").Value = 11239.64
$ws.Range("C5").Value = 11214.41

# Row 6: total
$ws.Range("A6").Value = "total"
$ws.Range("B6").Value = 33718.91
$ws.Range("C6").Value = 33643.23
